$wb = $excel.ActiveWorkbook

# --- Sheet "VDWorstCaseYellowPercentage" ---
# Add a (blank/space) device-name value in B6, bump the expected value in B8,
# and move the selected cell to C11.
$ws2 = $wb.Worksheets.Item("VDWorstCaseYellowPercentage")
$ws2.Activate()
$ws2.Range("B6").Value = " "
$ws2.Range("B8").Value = 47
$ws2.Range("C11").Select()

# --- Sheet "VtgDropYellowColorPercentage" ---
# Just move the selected cell to C25.
$ws3 = $wb.Worksheets.Item("VtgDropYellowColorPercentage")
$ws3.Activate()
$ws3.Range("C25").Select()

# --- Sheet "RedColorPercentage" ---
# This becomes the active/selected tab when the workbook is reopened.
$ws4 = $wb.Worksheets.Item("RedColorPercentage")
$ws4.Activate()
$ws4.Range("B4").Select()

$wb.Save()
